$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ingredients")

# Note: the shared-string table records new unique strings in the order
# they are first written, so write "halloumi cheese" before "green bananas"
# to reproduce the same sharedStrings.xml ordering as the target edit
# (halloumi cheese -> index 190, green bananas -> index 191).
$ws.Cells.Item($ws.Rows.Count, 50).Value = "halloumi cheese"
$ws.Cells.Item($ws.Rows.Count, 50).Value = "green bananas"
$ws.Cells.Item($ws.Rows.Count, 50).Value = ""

# --- Insert "green bananas" as a new row 81 ---
$ws.Rows.Item(81).Insert()
$ws.Cells.Item(81, 1).Value = "green bananas"
$ws.Cells.Item(81, 2).Value = "Vegetable"
$ws.Cells.Item(81, 3).Value = 88
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 23
$ws.Cells.Item(81, 6).Value = 0

# --- Insert "halloumi cheese" as a new row 94 ---
$ws.Rows.Item(94).Insert()
$ws.Cells.Item(94, 1).Value = "halloumi cheese"
$ws.Cells.Item(94, 2).Value = "Dairy"
$ws.Cells.Item(94, 3).Value = 316
$ws.Cells.Item(94, 4).Value = 20
$ws.Cells.Item(94, 5).Value = 2
$ws.Cells.Item(94, 6).Value = 25

# --- Update view state (pane / selection) to match the target ---
$ws.Application.ActiveWindow.ScrollRow = 143
$ws.Range("G172").Select()

# --- Keep the hidden AutoFilter defined name in sync with the new extent ---
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Ingredients!`$A`$2:`$G`$182"
